$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.683.96'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -5.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.209.08'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -6.88%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.94'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.85'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -10.35%  '
$ws.Range('E7').Value = '  -8.28%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.557'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -9.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.37'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -11.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.06'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0822'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -10.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.73'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -9.08%  '
$ws.Range('E14').Value = '  -4.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.860'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -12.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.544.67'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -6.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.07'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -7.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.201.03'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -7.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.576.16'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.74'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0956'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -10.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.36'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -12.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.91'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -11.53%  '
$ws.Range('E24').Value = '  -10.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '235.68'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -9.47%  '
$ws.Range('E26').Value = '  -8.63%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.00'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -10.39%  '
$ws.Range('E29').Value = '  -5.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.22'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -14.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.34'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -9.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0874'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -9.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '33.87'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -10.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '155.03'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -8.41%  '
$ws.Range('E35').Value = '  -6.56%  '
$ws.Range('E36').Value = '  +6.68%  '
$ws.Range('E37').Value = '  +11.47%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.122'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -6.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.42'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -6.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.102'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -12.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.70'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -6.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0324'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -8.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.867.48'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +6.61%  '
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '88.87'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -11.14%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.08'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -7.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.206'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -11.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.39'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '75.38'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -7.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '59.96'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -13.86%  '
$ws.Range('E51').Value = '  -6.61%  '
